# Update the multiplication operands/results shown in the practice-sheet
# table cells (three-digit number multiplied by one-digit number).
$d = $word.ActiveDocument
$d.Content.Find.Execute("420×5=", $true, $false, $false, $false, $false, $true, 1, $false, "837×8=", 2) | Out-Null
$d.Content.Find.Execute("126×8=", $true, $false, $false, $false, $false, $true, 1, $false, "238×7=", 2) | Out-Null
$d.Content.Find.Execute("741×7=", $true, $false, $false, $false, $false, $true, 1, $false, "343×5=", 2) | Out-Null
$d.Content.Find.Execute("788×8=", $true, $false, $false, $false, $false, $true, 1, $false, "459×5=", 2) | Out-Null
$d.Content.Find.Execute("722×3=", $true, $false, $false, $false, $false, $true, 1, $false, "516×9=", 2) | Out-Null
$d.Content.Find.Execute("675×3=", $true, $false, $false, $false, $false, $true, 1, $false, "776×8=", 2) | Out-Null
$d.Content.Find.Execute("537×5=", $true, $false, $false, $false, $false, $true, 1, $false, "398×9=", 2) | Out-Null
$d.Content.Find.Execute("572×9=", $true, $false, $false, $false, $false, $true, 1, $false, "468×2=", 2) | Out-Null
$d.Content.Find.Execute("330×6=", $true, $false, $false, $false, $false, $true, 1, $false, "538×5=", 2) | Out-Null
$d.Content.Find.Execute("207×5=", $true, $false, $false, $false, $false, $true, 1, $false, "192×8=", 2) | Out-Null
$d.Content.Find.Execute("141×7=", $true, $false, $false, $false, $false, $true, 1, $false, "707×7=", 2) | Out-Null
$d.Content.Find.Execute("237×9=", $true, $false, $false, $false, $false, $true, 1, $false, "269×9=", 2) | Out-Null
$d.Content.Find.Execute("832×5=", $true, $false, $false, $false, $false, $true, 1, $false, "472×3=", 2) | Out-Null
$d.Content.Find.Execute("326×9=", $true, $false, $false, $false, $false, $true, 1, $false, "881×9=", 2) | Out-Null
$d.Content.Find.Execute("779×7=", $true, $false, $false, $false, $false, $true, 1, $false, "151×5=", 2) | Out-Null
$d.Content.Find.Execute("223×3=", $true, $false, $false, $false, $false, $true, 1, $false, "569×4=", 2) | Out-Null
$d.Content.Find.Execute("511×6=", $true, $false, $false, $false, $false, $true, 1, $false, "174×2=", 2) | Out-Null
$d.Content.Find.Execute("691×3=", $true, $false, $false, $false, $false, $true, 1, $false, "195×5=", 2) | Out-Null
$d.Content.Find.Execute("545×9=", $true, $false, $false, $false, $false, $true, 1, $false, "519×3=", 2) | Out-Null
$d.Content.Find.Execute("394×8=", $true, $false, $false, $false, $false, $true, 1, $false, "872×4=", 2) | Out-Null
$d.Content.Find.Execute("181×3=", $true, $false, $false, $false, $false, $true, 1, $false, "193×4=", 2) | Out-Null
$d.Content.Find.Execute("585×5=", $true, $false, $false, $false, $false, $true, 1, $false, "965×3=", 2) | Out-Null
$d.Content.Find.Execute("318×8=", $true, $false, $false, $false, $false, $true, 1, $false, "358×4=", 2) | Out-Null
$d.Content.Find.Execute("854×3=", $true, $false, $false, $false, $false, $true, 1, $false, "966×5=", 2) | Out-Null
$d.Content.Find.Execute("482×9=", $true, $false, $false, $false, $false, $true, 1, $false, "889×3=", 2) | Out-Null
